# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# Refreshes the last 9 data rows (298-306) of the "Poland Ekstraklasa"
# sheet with updated match ids, teams and odds, mirroring a re-pull of
# the upstream odds feed (rows keep their position/index but the
# id/teams/odds attached to each position are reshuffled + updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, id, HomeTeam, AwayTeam, oddH_op, oddD_op, oddA_op, oddH, oddD, oddA,
# Ah, oddAHH, oddAHA, AhOU, oddAHOver, oddAHUnder
$data = @(
    @(298, '7093821', 'LKS Lodz', 'Stal Mielec', 2.5, 3.4, 2.5, 2.6, 3.4, 2.4, 0, 2, 1.85, 2.75, 1.875, 1.975),
    @(299, '7093820', 'Ruch Chorzow', 'Cracovia Krakow', 2.5, 3.4, 2.5, 2.625, 3.4, 2.4, 0, 2.05, 1.8, 2.5, 1.825, 2.025),
    @(300, '7090293', 'Radomiak Radom', 'Widzew Lodz', 2.2, 3.1, 3.1, 2.25, 3.1, 3, -0.25, 2, 1.85, 2.5, 1.8, 2.05),
    @(301, '7083189', 'Pogon Szczecin', 'Gornik Zabrze', 1.727, 4, 3.75, 1.8, 3.9, 3.6, -0.5, 1.825, 2.025, 3, 1.825, 2.025),
    @(302, '7083188', 'Legia Warsaw', 'Zaglebie Lubin', 1.5, 4, 5.5, 1.55, 3.9, 5, -1, 1.975, 1.875, 2.75, 1.925, 1.925),
    @(303, '7088350', 'Puszcza Niepolomice', 'Piast Gliwice', 3, 3.1, 2.3, 2.8, 3.1, 2.4, 0, 2.1, 1.775, 2.25, 1.875, 1.975),
    @(304, '7074364', 'Rakow Czestochowa', 'Slask Wroclaw', 2.5, 3.6, 2.4, 2.5, 3.6, 2.4, 0, 1.975, 1.875, 2.5, 1.85, 2),
    @(305, '7041338', 'Jagiellonia Bialystok', 'Warta Poznan', 1.444, 4.75, 5.25, 1.45, 4.75, 5.25, -1.25, 2.025, 1.825, 3, 1.925, 1.925),
    @(306, '7083187', 'Lech Poznan', 'Korona Kielce', 1.8, 3.8, 3.6, 1.9, 3.75, 3.4, -0.5, 1.95, 1.9, 2.75, 1.925, 1.925)
)

foreach ($r in $data) {
    $row = $r[0]

    # Column B ("id") is stored as text in the sheet (e.g. "7093820"), so
    # force a text format before assigning a numeric-looking string,
    # otherwise Excel would silently coerce it to a number cell. Clear the
    # formatting right after so the cell keeps its original (default) style
    # instead of picking up a new "text" number-format style.
    $idCell = $ws.Cells.Item($row, 2)
    $idCell.NumberFormat = "@"
    $idCell.Value = $r[1]
    $idCell.ClearFormats()

    $ws.Cells.Item($row, 5).Value = $r[2]    # HomeTeam
    $ws.Cells.Item($row, 6).Value = $r[3]    # AwayTeam

    $ws.Cells.Item($row, 10).Value = $r[4]   # oddH_op
    $ws.Cells.Item($row, 11).Value = $r[5]   # oddD_op
    $ws.Cells.Item($row, 12).Value = $r[6]   # oddA_op
    $ws.Cells.Item($row, 13).Value = $r[7]   # oddH
    $ws.Cells.Item($row, 14).Value = $r[8]   # oddD
    $ws.Cells.Item($row, 15).Value = $r[9]   # oddA
    $ws.Cells.Item($row, 16).Value = $r[10]  # Ah
    $ws.Cells.Item($row, 17).Value = $r[11]  # oddAHH
    $ws.Cells.Item($row, 18).Value = $r[12]  # oddAHA
    $ws.Cells.Item($row, 19).Value = $r[13]  # AhOU
    $ws.Cells.Item($row, 20).Value = $r[14]  # oddAHOver
    $ws.Cells.Item($row, 21).Value = $r[15]  # oddAHUnder
}
